# Update "Ziele Tabelle" worksheet: replace the old goals/requirements
# rows with the new project status rows, add one new row, fix up one
# cell's wrap formatting, and move the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Nr. 1) ---------------------------------------------------
$ws.Range("C3").Value = "Leistung"
$ws.Range("D3").Value = "Bug fix"
$ws.Range("E3").Value = "Alle Feher behoben, die wärend der Testphase im Kurse aufgekommen sind."
$ws.Range("F3").Value = "Muss"

# --- Row 4 (Nr. 2) ---------------------------------------------------
$ws.Range("C4").Value = "Leistung"
$ws.Range("D4").Value = "Bilder selber einfügen"
$ws.Range("E4").Value = "User kann selber bilder Hochladen für die Webseite."
$ws.Range("F4").Value = "Muss"

# --- Row 5 (Nr. 3) ---------------------------------------------------
$ws.Range("C5").Value = "Leistung"
$ws.Range("D5").Value = "User Funktionen erweitern"
$ws.Range("E5").Value = "Es sollen weitere Funktion wie Folgen, liken und Beachrichtigen hinzugefügt."
$ws.Range("F5").Value = "Muss"

# --- Row 6 (Nr. 4) ---------------------------------------------------
$ws.Range("C6").Value = "Leistung"
$ws.Range("D6").Value = "PDF überarbeiten "
$ws.Range("E6").Value = "PDF soll schöner werden. Der Name der PDF soll der Rezepte Titel sein."
$ws.Range("F6").Value = "Muss"

# --- Row 7 (Nr. 5) ---------------------------------------------------
$ws.Range("C7").Value = "Leistung "
$ws.Range("D7").Value = "Passwörter verschlüsseln"
$ws.Range("E7").Value = "Die Passwörter sollen verschlüsselt gesendet und gespeichert werden. "
$ws.Range("F7").Value = "Soll"

# --- Row 8 (Nr. 6) ---------------------------------------------------
$ws.Range("C8").Value = "Leistung"
$ws.Range("D8").Value = "Kategorien überarbeiten"
$ws.Range("E8").Value = "Die Kategorein sollen Nutzerfreundlich und mehr sinn ergeben."
$ws.Range("F8").Value = "Soll"
# D8 switches to the wrap-text variant of its fill style (matches E8).
$ws.Range("E8").Copy()
$ws.Range("D8").PasteSpecial(-4122)

# --- Row 9 (Nr. 7) ---------------------------------------------------
$ws.Range("C9").Value = "Leistung"
$ws.Range("D9").Value = "Kommentare"
$ws.Range("E9").Value = "Es soll möglich sein unter Rezepten zu Kommentieren."
$ws.Range("F9").Value = "Kann"

# --- Row 10 (Nr. 8, new row) ------------------------------------------
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "Leistung"
$ws.Range("D10").Value = "Rezepte bearbeiten"
$ws.Range("E10").Value = "Der User ist möglich seine erstelln Rezepte zu bearbeiten."
$ws.Range("F10").Value = "Kann"

# Row 10 had no formatting yet (it was a blank placeholder row) - copy
# the formatting used by the other even data rows (row 4's B/C/D/E/F
# pattern) so it matches the rest of the table.
$ws.Range("B4").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("F10").PasteSpecial(-4122)

# --- Selection moves from K4 to D3 ------------------------------------
$ws.Range("D3").Select() | Out-Null
